# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect every cell that changes, as "<ref>" = "<new value>" pairs.
$updates = [ordered]@{
    "D2" = '61.966.72'
    "E2" = '  -0.92%  '
    "D3" = '3.415.68'
    "E3" = '  -0.73%  '
    "E4" = '  -0.04%  '
    "D5" = '410.39'
    "E5" = '  +0.58%  '
    "D6" = '129.65'
    "E6" = '  -3.06%  '
    "D7" = '0.640'
    "E7" = '  +8.03%  '
    "E8" = '  -0.07%  '
    "D9" = '0.739'
    "E9" = '  +7.73%  '
    "E10" = '  +15.20%  '
    "D11" = '42.79'
    "E11" = '  +1.12%  '
    "D12" = '0.0000214'
    "E12" = '  +61.89%  '
    "D13" = '9.10'
    "E13" = '  +7.36%  '
    "E14" = '  -0.40%  '
    "D15" = '3.953.14'
    "E15" = '  -0.88%  '
    "D16" = '21.23'
    "E16" = '  +6.37%  '
    "D17" = '3.419.13'
    "E17" = '  -0.82%  '
    "D18" = '12.12'
    "E18" = '  +6.45%  '
    "E19" = '  +5.39%  '
    "D20" = '61.917.37'
    "E20" = '  -1.02%  '
    "D21" = '444.58'
    "E21" = '  +40.95%  '
    "D22" = '91.16'
    "E22" = '  +8.33%  '
    "D23" = '3.16'
    "E23" = '  -1.02%  '
    "D24" = '13.06'
    "E24" = '  +0.90%  '
    "D25" = '3.27'
    "E25" = '  +3.17%  '
    "D26" = '33.69'
    "E26" = '  +13.16%  '
    "D27" = '8.81'
    "E27" = '  +6.42%  '
    "E28" = '  +0.27%  '
    "D29" = '7.62'
    "E29" = '  +0.51%  '
    "D30" = '2.71'
    "E30" = '  -1.45%  '
    "D31" = '11.99'
    "E31" = '  +5.19%  '
    "E32" = '  -0.60%  '
    "D33" = '0.168'
    "E33" = '  -3.10%  '
    "D34" = '42.91'
    "D35" = '0.999'
    "E35" = '  -0.08%  '
    "D36" = '0.0501'
    "E36" = '  +3.00%  '
    "D37" = '53.46'
    "E37" = '  +3.98%  '
    "D38" = '0.999'
    "E38" = '  -0.03%  '
    "E39" = '  -1.20%  '
    "E40" = '  +7.59%  '
    "E41" = '  -1.39%  '
    "E42" = '  -1.73%  '
    "D43" = '141.31'
    "E43" = '  +2.54%  '
    "D44" = '4.21'
    "E44" = '  +4.52%  '
    "E45" = '  -0.67%  '
    "E46" = '  +7.91%  '
    "E47" = '  -1.19%  '
    "D48" = '22.24'
    "E48" = '  +3.60%  '
    "D49" = '3.765.28'
    "E49" = '  -0.53%  '
    "D50" = '2.105.78'
    "E50" = '  -1.15%  '
    "B51" = 'Cronos'
    "C51" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    "D51" = '0.133'
    "E51" = '  +13.77%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Keep the cell as plain text (matches the source data, which stores
    # prices/percentages as strings, e.g. "0.640", "9.10") so Excel does not
    # silently coerce numeric-looking text into a Number and drop trailing zeros.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
